# Remove the first data row (the "above 50 no upstream" entry).
# Excel shifts all remaining rows up by one, the sheet's dimension shrinks
# from A1:S6 to A1:S5, and on save the shared-strings table is rebuilt
# to drop the strings that were only referenced by the deleted row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Delete()
